$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.023.74"
$ws.Range("E2").Value = "  -14.06%  "
$ws.Range("D3").Value = "2.286.04"
$ws.Range("E3").Value = "  -21.15%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "439.71"
$ws.Range("E5").Value = "  -16.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "118.92"
$ws.Range("E6").Value = "  -16.73%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.462"
$ws.Range("E8").Value = "  -15.57%  "
$ws.Range("D9").Value = "2.298.82"
$ws.Range("E9").Value = "  -20.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.27"
$ws.Range("E10").Value = "  -11.72%  "
$ws.Range("E11").Value = "  -19.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.298"
$ws.Range("E12").Value = "  -16.61%  "
$ws.Range("E13").Value = "  -6.84%  "
$ws.Range("D14").Value = "52.115.11"
$ws.Range("E14").Value = "  -13.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.56"
$ws.Range("E15").Value = "  -17.30%  "
$ws.Range("E16").Value = "  -17.34%  "
$ws.Range("D17").Value = "2.313.72"
$ws.Range("E17").Value = "  -20.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.93"
$ws.Range("E18").Value = "  -20.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "294.54"
$ws.Range("E19").Value = "  -15.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.997"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.69"
$ws.Range("E21").Value = "  -24.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.64"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.07"
$ws.Range("E23").Value = "  -21.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "52.92"
$ws.Range("E24").Value = "  -18.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.363"
$ws.Range("E25").Value = "  -19.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.141"
$ws.Range("E26").Value = "  -20.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -12.30%  "
$ws.Range("D29").Value = "0.0₃0642"
$ws.Range("E29").Value = "  -24.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "142.58"
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "16.64"
$ws.Range("E31").Value = "  -14.55%  "
$ws.Range("E32").Value = "  -21.63%  "
$ws.Range("E33").Value = "  -16.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.803"
$ws.Range("E34").Value = "  -19.01%  "
$ws.Range("E35").Value = "  -22.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.980"
$ws.Range("E37").Value = "  -17.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.79"
$ws.Range("E38").Value = "  -15.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.11"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.554"
$ws.Range("E40").Value = "  -14.45%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0498"
$ws.Range("E41").Value = "  -13.93%  "
$ws.Range("E42").Value = "  -16.86%  "
$ws.Range("D43").Value = "1.896.55"
$ws.Range("E43").Value = "  -16.97%  "
$ws.Range("E44").Value = "  -21.32%  "
$ws.Range("E45").Value = "  -11.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0201"
$ws.Range("E46").Value = "  -14.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.05"
$ws.Range("E47").Value = "  -17.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.47"
$ws.Range("E48").Value = "  -23.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.61"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.42"
$ws.Range("E50").Value = "  -12.96%  "
$ws.Range("E51").Value = "  -19.09%  "
